$d = $word.ActiveDocument

# The template's greeting line reads "***Prenom*** ***Nom***" (a placeholder
# for the recipient's first/last name). Fix the misspelled "Prenom" to the
# correctly accented "Prénom". The replacement text spans across the
# original run boundaries, so Word folds the touched runs back together and
# drops the (now-stale) spell-check proofErr markers around the corrected
# word, exactly as it would if you retyped the word by hand.
$d.Content.Find.Execute("***Prenom***", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "***Prénom***", 2)
